$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowCells {
    param(
        [int]$Row1,
        [int]$Row2,
        [string[]]$Columns
    )

    # Read all values from both rows first, so writes to one row never
    # clobber data we still need to read from the other.
    $vals1 = @{}
    $vals2 = @{}
    foreach ($col in $Columns) {
        $vals1[$col] = $ws.Range("$col$Row1").Value()
        $vals2[$col] = $ws.Range("$col$Row2").Value()
    }

    foreach ($col in $Columns) {
        $ws.Range("$col$Row1").Value = $vals2[$col]
        $ws.Range("$col$Row2").Value = $vals1[$col]
    }
}

# Rows 3 and 4: full record swap (every populated column between the two rows).
# Note: Y/AA ("2026-01-21") are identical in both rows, so they are left out
# to avoid Excel's automatic text->date coercion on round-tripping the value.
$cols34 = @("A","B","D","E","F","G","H","I","P","Q","R","S","T","U","V","W","AC","AD","AE","AG","AT","AW","AX","AY")
Swap-RowCells 3 4 $cols34

# Rows 31 and 33: only the Id / coordinate columns differ, swap those
$cols3133 = @("A","Q","R")
Swap-RowCells 31 33 $cols3133

# Rows 32 and 34: full record swap (every populated column between the two rows)
$cols3234 = @("A","B","D","E","F","G","H","I","K","L","M","N","P","Q","R","S","T","U","V","W","AC","AD","AE","AG","AH","AJ","AK","AM","AO","AT","AW","AX","AY")
Swap-RowCells 32 34 $cols3234

Write-Host "Swap complete"
